$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.10618618800344
$ws.Range("C2").Value = 0.0485881989242378
$ws.Range("B3").Value = 19.92719279733022
$ws.Range("C3").Value = 0.006168607584611716
